$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Laptopdb")

$ws.Range("E2").Value = 16.0
$ws.Range("E4").Value = 99.0
$ws.Range("E5").Value = 99.0
$ws.Range("E8").Value = 99.0
